$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, M, N, P, R, S

# Column D (Fecha) - serial date numbers
$ws.Range("D2").Value = 44320
$ws.Range("D3").Value = 44362

# Column M (Volumen)
$ws.Range("M2").Value = 50
$ws.Range("M3").Value = 100

# Column N (Precio mínimo)
$ws.Range("N2").Value = 18000
$ws.Range("N3").Value = 19000

# Column P (Precio promedio ponderado)
$ws.Range("P2").Value = 18800
$ws.Range("P3").Value = 19500

# Column R (Origen)
$ws.Range("R2").Value = "Provincia de Limarí"
$ws.Range("R3").Value = "Provincia de Curicó"

# Column S (Precio $/Kg)
$ws.Range("S2").Value = 1044
$ws.Range("S3").Value = 1083
